# test_metadata.xlsx update: rename the "Site" header to "State",
# shrink the header row, and move the active selection to the
# now-edited header cell (H1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header H1 used to read "Site" -> now "State".
$ws.Range("H1").Value = "State"

# Header row height shrinks from 45.75 to 30.75.
$ws.Rows.Item(1).RowHeight = 30.75

# Selection moves from E6 to the edited header cell H1.
$ws.Range("H1").Select()
